{"js": "// Update the worksheet date and all 25 two-digit multiplication problems.\nconst replacements = [\n  [\"2025-11-14 Friday\", \"2025-11-15 Saturday\"],\n  [\"41\u00d791=\", \"31\u00d763=\"],\n  [\"60\u00d790=\", \"67\u00d743=\"],\n  [\"28\u00d728=\", \"37\u00d781=\"],\n  [\"28\u00d718=\", \"36\u00d788=\"],\n  [\"45\u00d723=\", \"82\u00d782=\"],\n  [\"71\u00d753=\", \"49\u00d745=\"],\n  [\"74\u00d756=\", \"30\u00d765=\"],\n  [\"13\u00d793=\", \"72\u00d785=\"],\n  [\"83\u00d717=\", \"60\u00d738=\"],\n  [\"24\u00d733=\", \"58\u00d773=\"],\n  [\"20\u00d785=\", \"60\u00d773=\"],\n  [\"41\u00d799=\", \"33\u00d758=\"],\n  [\"80\u00d715=\", \"82\u00d712=\"],\n  [\"70\u00d728=\", \"53\u00d745=\"],\n  [\"95\u00d724=\", \"48\u00d798=\"],\n  [\"45\u00d733=\", \"12\u00d739=\"],\n  [\"70\u00d782=\", \"25\u00d716=\"],\n  [\"48\u00d789=\", \"62\u00d747=\"],\n  [\"67\u00d736=\", \"76\u00d737=\"],\n  [\"63\u00d712=\", \"51\u00d724=\"],\n  [\"54\u00d713=\", \"26\u00d767=\"],\n  [\"31\u00d721=\", \"30\u00d772=\"],\n  [\"62\u00d726=\", \"68\u00d787=\"],\n  [\"67\u00d775=\", \"94\u00d784=\"],\n  [\"61\u00d790=\", \"81\u00d775=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and all 25 two-digit multiplication problems.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-11-14 Friday\", \"2025-11-15 Saturday\"),\n    @(\"41\u00d791=\", \"31\u00d763=\"),\n    @(\"60\u00d790=\", \"67\u00d743=\"),\n    @(\"28\u00d728=\", \"37\u00d781=\"),\n    @(\"28\u00d718=\", \"36\u00d788=\"),\n    @(\"45\u00d723=\", \"82\u00d782=\"),\n    @(\"71\u00d753=\", \"49\u00d745=\"),\n    @(\"74\u00d756=\", \"30\u00d765=\"),\n    @(\"13\u00d793=\", \"72\u00d785=\"),\n    @(\"83\u00d717=\", \"60\u00d738=\"),\n    @(\"24\u00d733=\", \"58\u00d773=\"),\n    @(\"20\u00d785=\", \"60\u00d773=\"),\n    @(\"41\u00d799=\", \"33\u00d758=\"),\n    @(\"80\u00d715=\", \"82\u00d712=\"),\n    @(\"70\u00d728=\", \"53\u00d745=\"),\n    @(\"95\u00d724=\", \"48\u00d798=\"),\n    @(\"45\u00d733=\", \"12\u00d739=\"),\n    @(\"70\u00d782=\", \"25\u00d716=\"),\n    @(\"48\u00d789=\", \"62\u00d747=\"),\n    @(\"67\u00d736=\", \"76\u00d737=\"),\n    @(\"63\u00d712=\", \"51\u00d724=\"),\n    @(\"54\u00d713=\", \"26\u00d767=\"),\n    @(\"31\u00d721=\", \"30\u00d772=\"),\n    @(\"62\u00d726=\", \"68\u00d787=\"),\n    @(\"67\u00d775=\", \"94\u00d784=\"),\n    @(\"61\u00d790=\", \"81\u00d775=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n"}
